$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue $ws "D2" "62.724.63"
$ws.Range("E2").Value = "  -0.47%  "

Set-TextValue $ws "D3" "3.042.98"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").Value = "  +0.20%  "

Set-TextValue $ws "D5" "536.25"
$ws.Range("E5").Value = "  -1.10%  "

Set-TextValue $ws "D6" "133.00"
$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("E7").Value = "  +0.08%  "

Set-TextValue $ws "D8" "3.041.06"
$ws.Range("E8").Value = "  -1.03%  "

Set-TextValue $ws "D9" "0.488"
$ws.Range("E9").Value = "  -0.29%  "

$ws.Range("E10").Value = "  -1.28%  "

Set-TextValue $ws "D11" "6.14"
$ws.Range("E11").Value = "  -2.22%  "

Set-TextValue $ws "D12" "0.446"
$ws.Range("E12").Value = "  -3.42%  "

$ws.Range("E13").Value = "  -3.05%  "

Set-TextValue $ws "D14" "33.78"
$ws.Range("E14").Value = "  -2.98%  "

Set-TextValue $ws "D15" "3.538.46"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D16" "0.112"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws "D17" "62.720.78"
$ws.Range("E17").Value = "  -0.37%  "

Set-TextValue $ws "D18" "3.048.30"
$ws.Range("E18").Value = "  -1.09%  "

Set-TextValue $ws "D19" "6.54"
$ws.Range("E19").Value = "  -1.66%  "

Set-TextValue $ws "D20" "464.39"
$ws.Range("E20").Value = "  -4.29%  "

$ws.Range("E21").Value = "  -1.10%  "

Set-TextValue $ws "D22" "0.684"
$ws.Range("E22").Value = "  -3.24%  "

$ws.Range("E23").Value = "  -4.59%  "

Set-TextValue $ws "D24" "77.58"
$ws.Range("E24").Value = "  -1.33%  "

Set-TextValue $ws "D25" "11.94"
$ws.Range("E25").Value = "  -1.79%  "

Set-TextValue $ws "D26" "0.997"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  -2.14%  "

Set-TextValue $ws "D28" "7.73"
$ws.Range("E28").Value = "  -5.91%  "

Set-TextValue $ws "D29" "0.999"
$ws.Range("E29").Value = "  +0.02%  "

Set-TextValue $ws "D30" "25.70"
$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D31" "1.85"
$ws.Range("E31").Value = "  -3.73%  "

$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D32" "1.13"
$ws.Range("E32").Value = "  +3.14%  "

Set-TextValue $ws "D33" "58.53"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("E34").Value = "  -6.92%  "

Set-TextValue $ws "D35" "5.33"
$ws.Range("E35").Value = "  +3.42%  "

$ws.Range("E36").Value = "  -2.70%  "

Set-TextValue $ws "D37" "463.79"
$ws.Range("E37").Value = "  -3.27%  "

Set-TextValue $ws "D38" "3.212.10"
$ws.Range("E38").Value = "  +2.29%  "

Set-TextValue $ws "D39" "0.0388"
$ws.Range("E39").Value = "  -0.79%  "

$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("E41").Value = "  -0.32%  "

Set-TextValue $ws "D42" "8.00"
$ws.Range("E42").Value = "  -1.10%  "

Set-TextValue $ws "D43" "2.51"
$ws.Range("E43").Value = "  -2.73%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("E45").Value = "  -2.80%  "

Set-TextValue $ws "D46" "121.72"
$ws.Range("E46").Value = "  +2.13%  "

Set-TextValue $ws "D47" "24.61"
$ws.Range("E47").Value = "  -0.72%  "

Set-TextValue $ws "D48" "0.107"
$ws.Range("E48").Value = "  -0.07%  "

Set-TextValue $ws "D49" "1.96"
$ws.Range("E49").Value = "  -4.18%  "

Set-TextValue $ws "D50" "0.0₃0509"
$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("E51").Value = "  +5.11%  "

